$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trial 2 (row 3) values: y_corrSteps, y_nrSteps, alienID
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Update selection to match saved view state (active cell E3)
$ws.Range("E3").Select()
